$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (border/number-format) for the new rows 1509-1559 by
# copying the existing blank-row format from row 1508 (the previous last row).
# This matches the style indices (s=3 for A-H, s=4 for I-J, s=8 once G/H hold
# True/False text) used throughout the rest of the sheet, instead of the bare
# column-default style new rows would otherwise pick up.
$ws.Range("A1508:J1508").Copy()
$ws.Range("A1509:J1559").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 1494: 2025-05-06
$ws.Range("A1494").Value = '2025-05-06'
$ws.Range("B1494").Value = '08:45'
$ws.Range("C1494").Value = '09:45'
$ws.Range("D1494").Value = '1h 00m'
$ws.Range("E1494").Value = '#studying'
$ws.Range("G1494").Value = "'False"
$ws.Range("H1494").Value = "'False"
$ws.Range("I1494").Formula = "=YEAR(A1494)"
$ws.Range("J1494").Formula = "=MONTH(A1494)"

# Row 1495: 2025-05-06
$ws.Range("A1495").Value = '2025-05-06'
$ws.Range("B1495").Value = '10:45'
$ws.Range("C1495").Value = '12:30'
$ws.Range("D1495").Value = '1h 45m'
$ws.Range("E1495").Value = '#studying'
$ws.Range("G1495").Value = "'False"
$ws.Range("H1495").Value = "'False"
$ws.Range("I1495").Formula = "=YEAR(A1495)"
$ws.Range("J1495").Formula = "=MONTH(A1495)"

# Row 1496: 2025-05-06
$ws.Range("A1496").Value = '2025-05-06'
$ws.Range("B1496").Value = '12:30'
$ws.Range("C1496").Value = '16:30'
$ws.Range("D1496").Value = '4h 00m'
$ws.Range("E1496").Value = '#maintenance'
$ws.Range("G1496").Value = "'False"
$ws.Range("H1496").Value = "'False"
$ws.Range("I1496").Formula = "=YEAR(A1496)"
$ws.Range("J1496").Formula = "=MONTH(A1496)"

# Row 1497: 2025-05-07
$ws.Range("A1497").Value = '2025-05-07'
$ws.Range("B1497").Value = '21:30'
$ws.Range("C1497").Value = '23:00'
$ws.Range("D1497").Value = '1h 30m'
$ws.Range("E1497").Value = '#maintenance'
$ws.Range("G1497").Value = "'False"
$ws.Range("H1497").Value = "'False"
$ws.Range("I1497").Formula = "=YEAR(A1497)"
$ws.Range("J1497").Formula = "=MONTH(A1497)"

# Row 1498: 2025-05-08
$ws.Range("A1498").Value = '2025-05-08'
$ws.Range("B1498").Value = '17:00'
$ws.Range("C1498").Value = '17:45'
$ws.Range("D1498").Value = '0h 45m'
$ws.Range("E1498").Value = '#python'
$ws.Range("F1498").Value = 'nwreadinglist v4.3.0'
$ws.Range("G1498").Value = "'True"
$ws.Range("H1498").Value = "'False"
$ws.Range("I1498").Formula = "=YEAR(A1498)"
$ws.Range("J1498").Formula = "=MONTH(A1498)"

# Row 1499: 2025-05-08
$ws.Range("A1499").Value = '2025-05-08'
$ws.Range("B1499").Value = '20:30'
$ws.Range("C1499").Value = '22:15'
$ws.Range("D1499").Value = '1h 45m'
$ws.Range("E1499").Value = '#python'
$ws.Range("F1499").Value = 'nwreadinglist v4.3.0'
$ws.Range("G1499").Value = "'True"
$ws.Range("H1499").Value = "'False"
$ws.Range("I1499").Formula = "=YEAR(A1499)"
$ws.Range("J1499").Formula = "=MONTH(A1499)"

# Row 1500: 2025-05-09
$ws.Range("A1500").Value = '2025-05-09'
$ws.Range("B1500").Value = '21:00'
$ws.Range("C1500").Value = '22:30'
$ws.Range("D1500").Value = '1h 30m'
$ws.Range("E1500").Value = '#python'
$ws.Range("F1500").Value = 'nwreadinglist v4.3.0'
$ws.Range("G1500").Value = "'True"
$ws.Range("H1500").Value = "'False"
$ws.Range("I1500").Formula = "=YEAR(A1500)"
$ws.Range("J1500").Formula = "=MONTH(A1500)"

# Row 1501: 2025-05-11
$ws.Range("A1501").Value = '2025-05-11'
$ws.Range("B1501").Value = '10:00'
$ws.Range("C1501").Value = '12:30'
$ws.Range("D1501").Value = '2h 30m'
$ws.Range("E1501").Value = '#python'
$ws.Range("F1501").Value = 'nwreadinglist v4.3.0'
$ws.Range("G1501").Value = "'True"
$ws.Range("H1501").Value = "'True"
$ws.Range("I1501").Formula = "=YEAR(A1501)"
$ws.Range("J1501").Formula = "=MONTH(A1501)"

# Row 1502: 2025-05-11
$ws.Range("A1502").Value = '2025-05-11'
$ws.Range("B1502").Value = '12:30'
$ws.Range("C1502").Value = '13:45'
$ws.Range("D1502").Value = '1h 15m'
$ws.Range("E1502").Value = '#maintenance'
$ws.Range("G1502").Value = "'False"
$ws.Range("H1502").Value = "'False"
$ws.Range("I1502").Formula = "=YEAR(A1502)"
$ws.Range("J1502").Formula = "=MONTH(A1502)"

# Row 1503: 2025-05-11
$ws.Range("A1503").Value = '2025-05-11'
$ws.Range("B1503").Value = '14:15'
$ws.Range("C1503").Value = '19:30'
$ws.Range("D1503").Value = '5h 15m'
$ws.Range("E1503").Value = '#maintenance'
$ws.Range("G1503").Value = "'False"
$ws.Range("H1503").Value = "'False"
$ws.Range("I1503").Formula = "=YEAR(A1503)"
$ws.Range("J1503").Formula = "=MONTH(A1503)"

# Row 1504: 2025-05-12
$ws.Range("A1504").Value = '2025-05-12'
$ws.Range("B1504").Value = '09:00'
$ws.Range("C1504").Value = '13:30'
$ws.Range("D1504").Value = '4h 30m'
$ws.Range("E1504").Value = '#maintenance'
$ws.Range("F1504").Value = 'Nuitka/pyinstaller. '
$ws.Range("G1504").Value = "'False"
$ws.Range("H1504").Value = "'False"
$ws.Range("I1504").Formula = "=YEAR(A1504)"
$ws.Range("J1504").Formula = "=MONTH(A1504)"

# Row 1505: 2025-05-12
$ws.Range("A1505").Value = '2025-05-12'
$ws.Range("B1505").Value = '14:00'
$ws.Range("C1505").Value = '17:30'
$ws.Range("D1505").Value = '3h 30m'
$ws.Range("E1505").Value = '#maintenance'
$ws.Range("F1505").Value = 'Nuitka/pyinstaller. '
$ws.Range("G1505").Value = "'False"
$ws.Range("H1505").Value = "'False"
$ws.Range("I1505").Formula = "=YEAR(A1505)"
$ws.Range("J1505").Formula = "=MONTH(A1505)"

# Row 1506: 2025-05-12
$ws.Range("A1506").Value = '2025-05-12'
$ws.Range("B1506").Value = '19:30'
$ws.Range("C1506").Value = '22:30'
$ws.Range("D1506").Value = '3h 00m'
$ws.Range("E1506").Value = '#maintenance'
$ws.Range("F1506").Value = 'Nuitka/pyinstaller. '
$ws.Range("G1506").Value = "'False"
$ws.Range("H1506").Value = "'False"
$ws.Range("I1506").Formula = "=YEAR(A1506)"
$ws.Range("J1506").Formula = "=MONTH(A1506)"

# Row 1507: 2025-05-13
$ws.Range("A1507").Value = '2025-05-13'
$ws.Range("B1507").Value = '10:30'
$ws.Range("C1507").Value = '13:30'
$ws.Range("D1507").Value = '3h 00m'
$ws.Range("E1507").Value = '#python'
$ws.Range("F1507").Value = 'nwtraderaanalytics v4.5.1'
$ws.Range("G1507").Value = "'True"
$ws.Range("H1507").Value = "'False"
$ws.Range("I1507").Formula = "=YEAR(A1507)"
$ws.Range("J1507").Formula = "=MONTH(A1507)"

# Row 1508: 2025-05-13
$ws.Range("A1508").Value = '2025-05-13'
$ws.Range("B1508").Value = '13:30'
$ws.Range("C1508").Value = '17:30'
$ws.Range("D1508").Value = '4h 00m'
$ws.Range("E1508").Value = '#python'
$ws.Range("F1508").Value = 'nwtraderaanalytics v4.5.1'
$ws.Range("G1508").Value = "'True"
$ws.Range("H1508").Value = "'False"
$ws.Range("I1508").Formula = "=YEAR(A1508)"
$ws.Range("J1508").Formula = "=MONTH(A1508)"

# Row 1509: 2025-05-14
$ws.Range("A1509").Value = '2025-05-14'
$ws.Range("B1509").Value = '17:45'
$ws.Range("C1509").Value = '19:45'
$ws.Range("D1509").Value = '2h 00m'
$ws.Range("E1509").Value = '#maintenance'
$ws.Range("F1509").Value = 'HackberryPi configuration.'
$ws.Range("G1509").Value = "'False"
$ws.Range("H1509").Value = "'False"
$ws.Range("I1509").Formula = "=YEAR(A1509)"
$ws.Range("J1509").Formula = "=MONTH(A1509)"

# Row 1510: 2025-05-15
$ws.Range("A1510").Value = '2025-05-15'
$ws.Range("B1510").Value = '17:00'
$ws.Range("C1510").Value = '17:45'
$ws.Range("D1510").Value = '0h 45m'
$ws.Range("E1510").Value = '#studying'
$ws.Range("G1510").Value = "'False"
$ws.Range("H1510").Value = "'False"
$ws.Range("I1510").Formula = "=YEAR(A1510)"
$ws.Range("J1510").Formula = "=MONTH(A1510)"

# Rows 1511/1512 are written with row 1512 first: matches the original
# authoring's shared-string insertion order (the "2025-05-17" string used
# by A1512 was registered before the "2025-05-16" string used by A1511).
# Row 1512: 2025-05-17
$ws.Range("A1512").Value = '2025-05-17'
$ws.Range("B1512").Value = '15:00'
$ws.Range("C1512").Value = '17:00'
$ws.Range("D1512").Value = '2h 00m'
$ws.Range("E1512").Value = '#maintenance'
$ws.Range("F1512").Value = 'HackberryPi configuration.'
$ws.Range("G1512").Value = "'False"
$ws.Range("H1512").Value = "'False"
$ws.Range("I1512").Formula = "=YEAR(A1512)"
$ws.Range("J1512").Formula = "=MONTH(A1512)"

# Row 1511: 2025-05-16
$ws.Range("A1511").Value = '2025-05-16'
$ws.Range("B1511").Value = '08:00'
$ws.Range("C1511").Value = '08:45'
$ws.Range("D1511").Value = '0h 45m'
$ws.Range("E1511").Value = '#studying'
$ws.Range("G1511").Value = "'False"
$ws.Range("H1511").Value = "'False"
$ws.Range("I1511").Formula = "=YEAR(A1511)"
$ws.Range("J1511").Formula = "=MONTH(A1511)"

# Row 1513: 2025-05-18
$ws.Range("A1513").Value = '2025-05-18'
$ws.Range("B1513").Value = '21:00'
$ws.Range("C1513").Value = '22:00'
$ws.Range("D1513").Value = '1h 00m'
$ws.Range("E1513").Value = '#python'
$ws.Range("F1513").Value = 'nwtraderaanalytics v4.5.1'
$ws.Range("G1513").Value = "'True"
$ws.Range("H1513").Value = "'True"
$ws.Range("I1513").Formula = "=YEAR(A1513)"
$ws.Range("J1513").Formula = "=MONTH(A1513)"

# Row 1514: 2025-05-18
$ws.Range("A1514").Value = '2025-05-18'
$ws.Range("B1514").Value = '22:45'
$ws.Range("C1514").Value = '23:45'
$ws.Range("D1514").Value = '1h 00m'
$ws.Range("E1514").Value = '#python'
$ws.Range("F1514").Value = 'nwdocstringchecking v1.0.0'
$ws.Range("G1514").Value = "'True"
$ws.Range("H1514").Value = "'False"
$ws.Range("I1514").Formula = "=YEAR(A1514)"
$ws.Range("J1514").Formula = "=MONTH(A1514)"

# Row 1515: 2025-05-19
$ws.Range("A1515").Value = '2025-05-19'
$ws.Range("B1515").Value = '11:45'
$ws.Range("C1515").Value = '17:15'
$ws.Range("D1515").Value = '5h 30m'
$ws.Range("E1515").Value = '#python'
$ws.Range("F1515").Value = 'nwdocstringchecking v1.0.0'
$ws.Range("G1515").Value = "'True"
$ws.Range("H1515").Value = "'True"
$ws.Range("I1515").Formula = "=YEAR(A1515)"
$ws.Range("J1515").Formula = "=MONTH(A1515)"

# Row 1516: 2025-05-19
$ws.Range("A1516").Value = '2025-05-19'
$ws.Range("B1516").Value = '17:15'
$ws.Range("C1516").Value = '18:15'
$ws.Range("D1516").Value = '1h 00m'
$ws.Range("E1516").Value = '#python'
$ws.Range("F1516").Value = 'nwcommitaverages v1.0.0'
$ws.Range("G1516").Value = "'True"
$ws.Range("H1516").Value = "'False"
$ws.Range("I1516").Formula = "=YEAR(A1516)"
$ws.Range("J1516").Formula = "=MONTH(A1516)"

# Row 1517: 2025-05-19
$ws.Range("A1517").Value = '2025-05-19'
$ws.Range("B1517").Value = '20:30'
$ws.Range("C1517").Value = '23:30'
$ws.Range("D1517").Value = '3h 00m'
$ws.Range("E1517").Value = '#python'
$ws.Range("F1517").Value = 'nwcommitaverages v1.0.0'
$ws.Range("G1517").Value = "'True"
$ws.Range("H1517").Value = "'False"
$ws.Range("I1517").Formula = "=YEAR(A1517)"
$ws.Range("J1517").Formula = "=MONTH(A1517)"

# Row 1518: 2025-05-20
$ws.Range("A1518").Value = '2025-05-20'
$ws.Range("B1518").Value = '10:00'
$ws.Range("C1518").Value = '16:00'
$ws.Range("D1518").Value = '6h 00m'
$ws.Range("E1518").Value = '#python'
$ws.Range("F1518").Value = 'nwcommitaverages v1.0.0'
$ws.Range("G1518").Value = "'True"
$ws.Range("H1518").Value = "'True"
$ws.Range("I1518").Formula = "=YEAR(A1518)"
$ws.Range("J1518").Formula = "=MONTH(A1518)"

# Row 1519: 2025-05-22
$ws.Range("A1519").Value = '2025-05-22'
$ws.Range("B1519").Value = '08:00'
$ws.Range("C1519").Value = '08:45'
$ws.Range("D1519").Value = '0h 45m'
$ws.Range("E1519").Value = '#studying'
$ws.Range("G1519").Value = "'False"
$ws.Range("H1519").Value = "'False"
$ws.Range("I1519").Formula = "=YEAR(A1519)"
$ws.Range("J1519").Formula = "=MONTH(A1519)"

# Row 1520: 2025-05-23
$ws.Range("A1520").Value = '2025-05-23'
$ws.Range("B1520").Value = '08:00'
$ws.Range("C1520").Value = '08:30'
$ws.Range("D1520").Value = '0h 30m'
$ws.Range("E1520").Value = '#studying'
$ws.Range("G1520").Value = "'False"
$ws.Range("H1520").Value = "'False"
$ws.Range("I1520").Formula = "=YEAR(A1520)"
$ws.Range("J1520").Formula = "=MONTH(A1520)"

# Row 1521: 2025-05-23
$ws.Range("A1521").Value = '2025-05-23'
$ws.Range("B1521").Value = '17:15'
$ws.Range("C1521").Value = '17:45'
$ws.Range("D1521").Value = '0h 30m'
$ws.Range("E1521").Value = '#studying'
$ws.Range("G1521").Value = "'False"
$ws.Range("H1521").Value = "'False"
$ws.Range("I1521").Formula = "=YEAR(A1521)"
$ws.Range("J1521").Formula = "=MONTH(A1521)"

$ws.Range("E1523").Select()
Write-Output "done"